# Apply cryptos list update (prices/volumes refreshed; a few coins reordered)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "69.886.98"
$ws.Cells.Item(2, 5).Value = "  -1.04%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.572.13"
$ws.Cells.Item(3, 5).Value = "  -1.88%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  +0.02%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'577.76"
$ws.Cells.Item(5, 5).Value = "  -3.17%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'188.13"
$ws.Cells.Item(6, 5).Value = "  -2.92%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -2.93%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "3.570.33"
$ws.Cells.Item(8, 5).Value = "  -0.90%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.999"
$ws.Cells.Item(9, 5).Value = "  -0.04%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'0.177"
$ws.Cells.Item(10, 5).Value = "  -3.25%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.658"
$ws.Cells.Item(11, 5).Value = "  -1.33%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'55.89"
$ws.Cells.Item(12, 5).Value = "  -4.10%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +0.16%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'9.66"
$ws.Cells.Item(14, 5).Value = "  -1.43%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "4.140.57"
$ws.Cells.Item(15, 5).Value = "  -1.71%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'19.92"
$ws.Cells.Item(16, 5).Value = "  +2.34%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "3.567.50"
$ws.Cells.Item(17, 5).Value = "  -1.62%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "69.770.66"
$ws.Cells.Item(18, 5).Value = "  -0.96%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'12.55"
$ws.Cells.Item(19, 5).Value = "  -1.20%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +0.10%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'1.04"
$ws.Cells.Item(21, 5).Value = "  -2.21%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'472.29"
$ws.Cells.Item(22, 5).Value = "  -5.02%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'19.17"
$ws.Cells.Item(23, 5).Value = "  +13.35%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -8.90%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'4.34"
$ws.Cells.Item(25, 5).Value = "  -3.25%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'88.43"
$ws.Cells.Item(26, 5).Value = "  -3.43%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'3.05"
$ws.Cells.Item(27, 5).Value = "  -2.74%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'10.98"
$ws.Cells.Item(28, 5).Value = "  -3.23%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'9.34"
$ws.Cells.Item(29, 5).Value = "  -0.92%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'32.10"
$ws.Cells.Item(30, 5).Value = "  -1.46%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +0.06%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +2.19%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'12.07"
$ws.Cells.Item(33, 5).Value = "  -1.58%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'65.68"
$ws.Cells.Item(34, 5).Value = "  -0.24%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'577.44"
$ws.Cells.Item(35, 5).Value = "  -6.78%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'38.48"
$ws.Cells.Item(36, 5).Value = "  +0.39%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'1.00"
$ws.Cells.Item(37, 5).Value = "  -0.12%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "0.0₃0799"
$ws.Cells.Item(38, 5).Value = "  -4.43%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -2.24%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.140"
$ws.Cells.Item(40, 5).Value = "  -6.08%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "dogwifhat"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(41, 4).Value = "'3.22"
$ws.Cells.Item(41, 5).Value = "  +14.58%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Stacks"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(42, 4).Value = "'3.49"
$ws.Cells.Item(42, 5).Value = "  -6.52%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "3.226.46"
$ws.Cells.Item(43, 5).Value = "  -3.34%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "Fetch.AI"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(44, 4).Value = "'2.85"
$ws.Cells.Item(44, 5).Value = "  +5.37%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.0441"
$ws.Cells.Item(46, 5).Value = "  -1.60%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'9.49"
$ws.Cells.Item(47, 5).Value = "  +4.35%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +1.02%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -0.66%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.999"
$ws.Cells.Item(50, 5).Value = "  +0.00%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -3.99%  "
